# updates to 5th, 4th, 3rd edition for radial scars
#
# This script:
#  1. Inserts two new rows into the "Edited_Tumors" sheet (right after the
#     existing "radial scar and complex sclerosing lesion" entry) for the
#     new, separated tumor names "radial scar " and "complex sclerosing
#     lesion" -- this also appends the two new strings to the shared
#     string table.
#  2. Consolidates the three per-column ("Yes" highlighted in pink)
#     conditional formatting rules on the "Generated" sheet (columns B, C,
#     D) into a single rule covering B:D.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Edited_Tumors: split "radial scar and complex sclerosing lesion"
#    into two new rows right after it.
# ---------------------------------------------------------------------
$editedTumors = $wb.Worksheets.Item("Edited_Tumors")

$editedTumors.Range("A2081:A2082").Insert()
$editedTumors.Range("A2081").Value = "radial scar "
$editedTumors.Range("A2082").Value = "complex sclerosing lesion"

# best-effort view-state update to match where the author was working
$editedTumors.Activate()
[void]$editedTumors.Range("A2072").Select()
$excel.ActiveWindow.Zoom = 171

# ---------------------------------------------------------------------
# 2. Generated: merge the B/C/D "contains Yes" conditional formats into
#    one rule applied across B1:D1048576.
# ---------------------------------------------------------------------
$generated = $wb.Worksheets.Item("Generated")

$rngB = $generated.Range("B1:B1048576")
$rngC = $generated.Range("C1:C1048576")
$rngD = $generated.Range("D1:D1048576")

$fcD = $rngD.FormatConditions.Item(1)
$fcD.ModifyAppliesToRange($generated.Range("B1:D1048576"))
$fcD.Formula1 = '=NOT(ISERROR(SEARCH("Yes",B1)))'
$fcD.Priority = 1

$rngB.FormatConditions.Delete()
$rngC.FormatConditions.Delete()
